# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
# The underlying worker database changed, so the "Estado de Cuenta" table
# (rows 16-24, columns C:G = Doc Trabajador, Nombre Trabajador, Periodo Mora,
# Valor Mora, Salario Basico) is refreshed/re-sorted with the latest values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value2 = "79498025"
$ws.Range("D16").Value2 = "JULIO ANDRES APRAEZ ESPAÑA"
$ws.Range("E16").Value2 = "2210"
$ws.Range("F16").Value2 = 82666
$ws.Range("G16").Value2 = 3600000
$ws.Range("C17").Value2 = "73168115"
$ws.Range("D17").Value2 = "LUIS ENRIQUE ABELLO MENDOZA"
$ws.Range("E17").Value2 = "2112"
$ws.Range("F17").Value2 = 45760
$ws.Range("G17").Value2 = 1000000
$ws.Range("C18").Value2 = "73182225"
$ws.Range("D18").Value2 = "POLICARPO MARTINEZ GARCIA"
$ws.Range("E18").Value2 = "2112"
$ws.Range("F18").Value2 = 36341
$ws.Range("G18").Value2 = 1480000
$ws.Range("C19").Value2 = "73182225"
$ws.Range("D19").Value2 = "POLICARPO MARTINEZ GARCIA"
$ws.Range("E19").Value2 = "2111"
$ws.Range("F19").Value2 = 6057
$ws.Range("G19").Value2 = 1480000
$ws.Range("C20").Value2 = "9149836"
$ws.Range("D20").Value2 = "NELSON YAIR ARIZA ANZOATEGUI"
$ws.Range("E20").Value2 = "2504"
$ws.Range("F20").Value2 = 80000
$ws.Range("G20").Value2 = 2000000
$ws.Range("C21").Value2 = "9149836"
$ws.Range("D21").Value2 = "NELSON YAIR ARIZA ANZOATEGUI"
$ws.Range("E21").Value2 = "2412"
$ws.Range("F21").Value2 = 80000
$ws.Range("G21").Value2 = 2000000
$ws.Range("C22").Value2 = "1047402132"
$ws.Range("D22").Value2 = "ELIANA ANDREA LORDUY TOBIO"
$ws.Range("E22").Value2 = "2203"
$ws.Range("F22").Value2 = 44000
$ws.Range("G22").Value2 = 1100000
$ws.Range("C23").Value2 = "1047402132"
$ws.Range("D23").Value2 = "ELIANA ANDREA LORDUY TOBIO"
$ws.Range("E23").Value2 = "2202"
$ws.Range("F23").Value2 = 44000
$ws.Range("G23").Value2 = 1100000
$ws.Range("C24").Value2 = "77094911"
$ws.Range("D24").Value2 = "JUSUE JOSE ONATE ROQUE"
$ws.Range("E24").Value2 = "2111"
$ws.Range("F24").Value2 = 6057
$ws.Range("G24").Value2 = 1400000
